# Add two new columns, I (I0) and J (IF), to the data table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: copy formatting from the existing header style (column H)
# so the new header cells pick up the same bold/border/center style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-34
$data = @(
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(1, 3),
    @(6, 7),
    @(10, 10),
    @(6, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 9),
    @(6, 7),
    @(6, 8),
    @(4, 6),
    @(6, 7),
    @(4, 5),
    @(5, 6),
    @(9, 9),
    @(3, 6),
    @(4, 5),
    @(4, 6),
    @(7, 7),
    @(8, 8),
    @(2, 4),
    @(8, 8),
    @(6, 6),
    @(8, 9),
    @(9, 9),
    @(5, 6),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
